$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Append four more measurement readings below the existing data (A2:A101)
$values = @(5, 3, 4, 5)
$startRow = 102

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Scroll the view down so the newly entered rows are visible, and leave the
# selection on the next empty row below the appended data (A106) - mirrors
# what Excel does after typing values down a column.
$excel.ActiveWindow.ScrollRow = 97
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A106").Select()
